$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new table row above row 24 (shifts rows 24:132 down to 25:133) ---
$ws.Rows("24:24").Insert()

# Re-attach the Table1 range so it covers the newly inserted row (A8:K133)
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A8:K133"))

# --- Fix up the brand-new row 24 (year separator row, like row 10's "2023") ---
# Copy the normal data-row formatting (row 23) into the new row first
$ws.Range("A23:K23").Copy()
$ws.Range("A24:K24").PasteSpecial(-4122)
# Restore the calculated-column formula in column G that PasteSpecial(formats) doesn't carry
$ws.Range("G24").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
# Now apply the special "year heading" formatting (from A10, which holds "2023") to A24 only
$ws.Range("A10").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A24").Value = "'2024"

# The table's auto-fill of the new last row (133) can mangle its formula; restore it
$ws.Range("G133").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- Fill in the previously-missing monthly EARNED entries ---
$ws.Range("C20").Value = 1.25
$ws.Range("C22").Value = 1.25
$ws.Range("C23").Value = 1.25

# --- New leave entry for 1/1/2024 (now row 25, shifted down from row 24) ---
$ws.Range("B25").Value = "VL(1-0-0)"
$ws.Range("D25").Value = 1
# K25 needs the date-stamp formatting used elsewhere in the REMARKS column (e.g. K10)
$ws.Range("K10").Copy()
$ws.Range("K25").PasteSpecial(-4122)
$ws.Range("K25").Value = 45300

$ws.Activate()
$ws.Range("H30").Select()

$excel.CutCopyMode = 0
